# Add group-header support to the Tom & Jerry table:
#  - a new "Category" column (D) and a new sub-header row (3)
#  - B2 becomes "Generations" (was "Age"); old "Age" value moves down to B3
#  - C2 becomes "Generations" as well; a brand new "Birth" sub-header lands in C3
#  - D2/D3 merge to host the (relocated) "Category" super-header
#  - A2/A3 merge so "Name" still spans both header rows
#  - the red-filled "Generations" cell (column B) switches from bold/black
#    text to non-bold/white text, and gains center/center alignment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Relabel the existing header row first (while the sheet is still
#    3 columns / 4 rows) so later structural inserts can just inherit
#    the already-correct formatting from their neighbours.
# ------------------------------------------------------------------
$ws.Range("B2").Value2 = "Generations"
$ws.Range("C2").Value2 = "Generations"

$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("C2").PasteSpecial(-4122)

$c = $ws.Range("B2")
$c.Font.Bold = $false
$c.Font.Color = 16777215
$c.HorizontalAlignment = -4108
$c.VerticalAlignment = -4108

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Structural inserts: a new sub-header row first, THEN the new
#    "Category" column -- in this order every newly created cell
#    inherits formatting that is already final, instead of stale
#    pre-edit formatting.
# ------------------------------------------------------------------
$ws.Rows("3").Insert()
$ws.Rows("3").RowHeight = 30
$ws.Columns("D").Insert()

# ------------------------------------------------------------------
# 3. Re-merge the title, and merge the two rowspan header cells
# ------------------------------------------------------------------
$ws.Range("A1:C1").UnMerge()
$ws.Range("A1:D1").Merge()
$ws.Range("A2:A3").Merge()
$ws.Range("D2:D3").Merge()

# ------------------------------------------------------------------
# 4. Remaining cell values (new column/row content)
# ------------------------------------------------------------------
$ws.Range("D2").Value2 = "Category"
$ws.Range("B3").Value2 = "Age"
$ws.Range("C3").Value2 = "Birth"

# ------------------------------------------------------------------
# 5. Column widths (A:D)
# ------------------------------------------------------------------
$ws.Columns("A:D").ColumnWidth = 19.9
